# "raw data + loader"
#
# Rename the two worksheets to reflect their actual roles, make the
# "label" sheet the active/selected tab, and re-prioritise the colour
# scale so the custom (theme-coloured) rule evaluates first.

$wb = $excel.ActiveWorkbook

$dataSheet  = $wb.Worksheets.Item(1)
$labelSheet = $wb.Worksheets.Item(2)

# --- sheet renames --------------------------------------------------
$dataSheet.Name  = "data"
$labelSheet.Name = "label"

# --- move the active tab from "data" to "label" ----------------------
# Activating the sheet updates both the workbook's activeTab and moves
# tabSelected="1" onto this sheet's own sheetView.
$labelSheet.Activate()

# --- reorder the B-column colour-scale conditional formatting --------
# Existing rules (in file order) are:
#   1) priority 2 - red/yellow/green percentile scale
#   2) priority 3 - green/cream 2-colour scale
#   3) priority 1 - custom theme-coloured percentile scale
# The custom (priority 1) rule should become the FIRST rule in the
# sheet's conditional-formatting list. Delete the other two rules (the
# theme-coloured rule's own position/definition is left untouched) and
# re-add them after it, restoring the original colours.
$cfRange = $dataSheet.Range("B1:B1048576")
$fcs = $cfRange.FormatConditions

$fcs.Item(1).Delete()
$fcs.Item(1).Delete()

$threeColor = $fcs.AddColorScale(3)
$twoColor = $fcs.AddColorScale(2)
